$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60 holds the new data point. Column A is a date-formatted label that
# must stay literal text (matching the rest of column A), so we force text
# entry with a leading apostrophe and then reset the cell style back to
# "Normal" so no quote-prefix style lingers on the cell.
$ws.Cells.Item(60, 1).Value = "'2025-10-14"
$ws.Cells.Item(60, 1).Style = "Normal"

$ws.Cells.Item(60, 2).Value = 53.95000076293945
$ws.Cells.Item(60, 3).Value = 395.4500122070312
$ws.Cells.Item(60, 4).Value = 347.75
